# Expenses Details.xlsx update
#  - Rename Sheet1 -> Expenses
#  - Add new sheet "Internal" after it, with individual expense settlements
#  - Update a handful of cells/formulas on the Expenses sheet

$wb = $excel.ActiveWorkbook
$expenses = $wb.Worksheets.Item(1)
$expenses.Name = "Expenses"

# ---------------------------------------------------------------------------
# Expenses sheet data updates
# ---------------------------------------------------------------------------

# Row 7: petrol-expenses row no longer attributed to Parag -- zeroed out
$expenses.Range("E7").Value = 0
$expenses.Range("F7").Value = 0

# Row 8: new invoice entry paid by Soby, with a link to the bill + updated total received
$expenses.Range("F8").Value = "Soby"
$expenses.Range("G8").Value = "Invoice"
$expenses.Hyperlinks.Add($expenses.Range("G8"), "https://github.com/Akshay1595/Lithium-ion_battery_as_service/blob/master/Docs/Bills/Bill_6.JPG") | Out-Null
$expenses.Range("I8").Value = 29653

# Header cells for the two summary columns on the right
$expenses.Range("L1").Value = "Individual"
$expenses.Range("M1").Value = "Total Expenses"

# Balance formula now computes Received - Spent instead of Spent - Received
$expenses.Range("H11").Formula = "=E11-E10"

# Selection cursor left where the author last clicked
$expenses.Range("L6").Select() | Out-Null

# Column widths - best effort autosize approximations for the updated columns
$expenses.Columns.Item(3).ColumnWidth = 24.833333333333336
$expenses.Columns.Item(12).ColumnWidth = 9
$expenses.Columns.Item(13).ColumnWidth = 13.5

# ---------------------------------------------------------------------------
# New "Internal" sheet: individual-to-individual settlements
# ---------------------------------------------------------------------------

$internal = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $expenses)
$internal.Name = "Internal"

$internal.Range("A1").Value = "Sr No."
$internal.Range("B1").Value = "Amount "
$internal.Range("C1").Value = "From Individual"
$internal.Range("D1").Value = "To Idividual"
$internal.Range("E1").Value = "Date"
$internal.Range("F1").Value = "Remark"

$internal.Range("A2").Value = 1
$internal.Range("B2").Value = 20000
$internal.Range("C2").Value = "Soby"
$internal.Range("D2").Value = "Parag"
$internal.Range("F2").Value = "Paid as a kick start amount"

$internal.Range("A3").Value = 2
$internal.Range("B3").Value = 5000
$internal.Range("C3").Value = "Soby"
$internal.Range("D3").Value = "Parag"
$internal.Range("F3").Value = "Paid Procurement Charges for initial components and testing "

$internal.Range("A4").Value = 3
$internal.Range("B4").Value = 20000
$internal.Range("C4").Value = "Parag"
$internal.Range("D4").Value = "Akshay"
$internal.Range("F4").Value = "Paid as a share of Engineering Cost"

# Header styling: bold, yellow fill, thin border - matches the Expenses header look
$hdr = $internal.Range("A1:F1")
$hdr.Font.Bold = $true
$hdr.Interior.Color = $expenses.Range("A1").Interior.Color
$hdr.Borders.Color = $expenses.Range("A1").Borders.Item(1).Color
$hdr.Borders.LineStyle = $expenses.Range("A1").Borders.Item(1).LineStyle
$hdr.Borders.Weight = $expenses.Range("A1").Borders.Item(1).Weight

# Body rows 2-23: thin border around every cell, matching the Expenses sheet style
$body = $internal.Range("A2:F23")
$body.Borders.Color = $expenses.Range("A2").Borders.Item(1).Color
$body.Borders.LineStyle = $expenses.Range("A2").Borders.Item(1).LineStyle
$body.Borders.Weight = $expenses.Range("A2").Borders.Item(1).Weight

# Column widths observed in the target sheet
$internal.Columns.Item(3).ColumnWidth = 14.166666666666666
$internal.Columns.Item(4).ColumnWidth = 10.5
$internal.Columns.Item(6).ColumnWidth = 55.666666666666664

$internal.Range("A1").Select() | Out-Null
